$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 ("New York -- New York") re-run failed with a GitHub API rate-limit error.
# Clear the previously-fetched values back to empty (inline string) cells,
# flip "Pct Includes Hispanic Black" (J4) from TRUE to FALSE, and update
# the status message in O4.

$ws.Range("B4:H4").Clear()
$ws.Range("J4").Value = $false
$ws.Range("K4:L4").Clear()

$ws.Range("O4").Value = "An error occurred. ... RateLimitExceededException(403, {'message': ""API rate limit exceeded for 132.145.200.60. (But here's the good news: Authenticated requests get a higher rate limit. Check out the documentation for more details.)"", 'documentation_url': 'https://developer.github.com/v3/#rate-limiting'})"
